# Fruta / hortaliza, semanal
# Insert a new weekly record as row 66 on the "Naranja" sheet, pushing the
# existing rows 66-90 down to 67-91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 66 (copies formatting/style from
# the row above, matching the style duplication seen in the target file).
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly record.
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 44755
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100102
$ws.Range("H66").Value = "Cítricos"
$ws.Range("I66").Value = 100102005
$ws.Range("J66").Value = "Naranja"
$ws.Range("K66").Value = "Fukumoto"
$ws.Range("L66").Value = "Tercera"
$ws.Range("M66").Value = 300
$ws.Range("N66").Value = 600
$ws.Range("O66").Value = 650
$ws.Range("P66").Value = 625
$ws.Range("Q66").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R66").Value = "Región de O'Higgins"
$ws.Range("S66").Value = 625
$ws.Range("T66").Value = 1

# Ensure the date cell keeps the date/time number format used by column D.
$ws.Range("D66").NumberFormat = $ws.Range("D67").NumberFormat
